$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $val) {
    $r = $ws.Range($rangeAddr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue 'D2' '70.361.88'
Set-TextValue 'E2' '  +0.71%  '
Set-TextValue 'D3' '3.620.44'
Set-TextValue 'E3' '  +2.70%  '
Set-TextValue 'E4' '  +0.12%  '
Set-TextValue 'D5' '601.80'
Set-TextValue 'E5' '  -0.74%  '
Set-TextValue 'D6' '196.64'
Set-TextValue 'E6' '  -0.30%  '
Set-TextValue 'E7' '  -0.80%  '
Set-TextValue 'E8' '  +0.06%  '
Set-TextValue 'E9' '  +6.96%  '
Set-TextValue 'D10' '0.646'
Set-TextValue 'E10' '  -0.35%  '
Set-TextValue 'D11' '53.24'
Set-TextValue 'E11' '  -1.05%  '
Set-TextValue 'E12' '  +0.89%  '
Set-TextValue 'D13' '9.55'
Set-TextValue 'E13' '  +0.32%  '
Set-TextValue 'D14' '4.195.39'
Set-TextValue 'E14' '  +2.72%  '
Set-TextValue 'D15' '605.92'
Set-TextValue 'E15' '  +1.43%  '
Set-TextValue 'E16' '  +1.28%  '
Set-TextValue 'D17' '70.446.91'
Set-TextValue 'E17' '  +0.61%  '
Set-TextValue 'D18' '3.623.62'
Set-TextValue 'E18' '  +2.72%  '
Set-TextValue 'D19' '19.03'
Set-TextValue 'E20' '  +1.42%  '
Set-TextValue 'D21' '0.998'
Set-TextValue 'E21' '  +0.61%  '
Set-TextValue 'D22' '18.16'
Set-TextValue 'E22' '  -1.24%  '
Set-TextValue 'D23' '5.20'
Set-TextValue 'E23' '  -1.68%  '
Set-TextValue 'D24' '103.12'
Set-TextValue 'E24' '  +1.25%  '
Set-TextValue 'D25' '4.60'
Set-TextValue 'E25' '  -0.46%  '
Set-TextValue 'E26' '  -7.10%  '
Set-TextValue 'D27' '10.62'
Set-TextValue 'E27' '  -2.56%  '
Set-TextValue 'D28' '9.69'
Set-TextValue 'E28' '  +0.88%  '
Set-TextValue 'E29' '  +1.34%  '
Set-TextValue 'D30' '4.67'
Set-TextValue 'E30' '  +7.52%  '
Set-TextValue 'E31' '  +2.82%  '
Set-TextValue 'D32' '12.27'
Set-TextValue 'E32' '  -1.48%  '
Set-TextValue 'E33' '  +1.04%  '
Set-TextValue 'D34' '63.30'
Set-TextValue 'E34' '  +0.24%  '
Set-TextValue 'D35' '0.0₃0889'
Set-TextValue 'E35' '  +3.55%  '
Set-TextValue 'D36' '3.923.00'
Set-TextValue 'E36' '  +5.72%  '
Set-TextValue 'B37' 'Bittensor'
Set-TextValue 'C37' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D37' '521.20'
Set-TextValue 'E37' '  +6.63%  '
Set-TextValue 'B38' 'Dai'
Set-TextValue 'C38' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D38' '0.998'
Set-TextValue 'E38' '  -0.14%  '
Set-TextValue 'E39' '  -0.24%  '
Set-TextValue 'D40' '36.70'
Set-TextValue 'E40' '  +0.14%  '
Set-TextValue 'E41' '  -1.16%  '
Set-TextValue 'E42' '  -2.56%  '
Set-TextValue 'D43' '0.136'
Set-TextValue 'E43' '  +2.61%  '
Set-TextValue 'D44' '0.0460'
Set-TextValue 'E44' '  +1.45%  '
Set-TextValue 'D45' '3.48'
Set-TextValue 'E45' '  +5.28%  '
Set-TextValue 'E46' '  +2.18%  '
Set-TextValue 'E47' '  -0.21%  '
Set-TextValue 'E48' '  +0.43%  '
Set-TextValue 'E49' '  -0.19%  '
Set-TextValue 'D50' '0.000249'
Set-TextValue 'E50' '  -0.41%  '
Set-TextValue 'E51' '  +0.63%  '
